$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell in this sheet is stored as literal text (prices use "."
# as a thousands separator, e.g. "65.836.18", and the 1h-volume column
# is a padded "  +1.23%  " string), so plain values must round-trip as
# text even when they look like a simple decimal (e.g. "583.03"). Excel
# auto-converts a General-formatted cell to a Number when the new text
# parses as one, so pre-mark just those cells as Text, write the value,
# then restore the Normal style so no stray formatting is left behind.
$textCells = @('D5', 'D6', 'D7', 'D10', 'D11', 'D15', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D29', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D51')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin prices / 1h volume changes (and the swapped
# Litecoin / Dai rows) from the latest GitHub Actions refresh.
$ws.Range('D2').Value2 = '65.836.18'
$ws.Range('E2').Value2 = '  -0.86%  '
$ws.Range('D3').Value2 = '3.279.62'
$ws.Range('E3').Value2 = '  -0.67%  '
$ws.Range('E4').Value2 = '  -0.02%  '
$ws.Range('D5').Value2 = '583.03'
$ws.Range('E5').Value2 = '  +2.11%  '
$ws.Range('D6').Value2 = '179.75'
$ws.Range('E6').Value2 = '  -1.06%  '
$ws.Range('D7').Value2 = '0.641'
$ws.Range('E7').Value2 = '  +7.32%  '
$ws.Range('E8').Value2 = '  -0.04%  '
$ws.Range('E9').Value2 = '  -3.40%  '
$ws.Range('D10').Value2 = '6.75'
$ws.Range('E10').Value2 = '  +1.99%  '
$ws.Range('D11').Value2 = '0.401'
$ws.Range('E11').Value2 = '  -0.26%  '
$ws.Range('D12').Value2 = '3.847.48'
$ws.Range('E12').Value2 = '  -0.73%  '
$ws.Range('E13').Value2 = '  -4.55%  '
$ws.Range('D14').Value2 = '65.934.50'
$ws.Range('E14').Value2 = '  -0.83%  '
$ws.Range('D15').Value2 = '26.02'
$ws.Range('E15').Value2 = '  -3.92%  '
$ws.Range('D16').Value2 = '3.263.01'
$ws.Range('E16').Value2 = '  -0.91%  '
$ws.Range('E17').Value2 = '  -2.76%  '
$ws.Range('D18').Value2 = '425.77'
$ws.Range('E18').Value2 = '  -0.75%  '
$ws.Range('D19').Value2 = '13.18'
$ws.Range('E19').Value2 = '  -3.25%  '
$ws.Range('D20').Value2 = '5.49'
$ws.Range('E20').Value2 = '  -3.08%  '
$ws.Range('D21').Value2 = '7.36'
$ws.Range('E21').Value2 = '  -3.08%  '
$ws.Range('B22').Value2 = 'Litecoin'
$ws.Range('C22').Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').Value2 = '71.59'
$ws.Range('E22').Value2 = '  -2.56%  '
$ws.Range('B23').Value2 = 'Dai'
$ws.Range('C23').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value2 = '1.00'
$ws.Range('E23').Value2 = '  +0.20%  '
$ws.Range('D24').Value2 = '5.68'
$ws.Range('E24').Value2 = '  +0.47%  '
$ws.Range('D25').Value2 = '3.419.65'
$ws.Range('E25').Value2 = '  -0.90%  '
$ws.Range('D26').Value2 = '0.508'
$ws.Range('E26').Value2 = '  -1.36%  '
$ws.Range('D27').Value2 = '0.197'
$ws.Range('E27').Value2 = '  +1.92%  '
$ws.Range('E28').Value2 = '  -4.58%  '
$ws.Range('D29').Value2 = '8.88'
$ws.Range('E29').Value2 = '  -1.39%  '
$ws.Range('E31').Value2 = '  -0.68%  '
$ws.Range('D32').Value2 = '22.20'
$ws.Range('E32').Value2 = '  -2.28%  '
$ws.Range('D33').Value2 = '1.00'
$ws.Range('E33').Value2 = '  +0.09%  '
$ws.Range('D34').Value2 = '5.14'
$ws.Range('E34').Value2 = '  -3.03%  '
$ws.Range('D35').Value2 = '6.57'
$ws.Range('E35').Value2 = '  -2.77%  '
$ws.Range('E36').Value2 = '  -3.76%  '
$ws.Range('D37').Value2 = '158.89'
$ws.Range('E37').Value2 = '  -0.55%  '
$ws.Range('D38').Value2 = '1.43'
$ws.Range('E38').Value2 = '  -4.46%  '
$ws.Range('E39').Value2 = '  -2.61%  '
$ws.Range('D40').Value2 = '26.31'
$ws.Range('E40').Value2 = '  -2.85%  '
$ws.Range('D41').Value2 = '2.794.09'
$ws.Range('E41').Value2 = '  -0.01%  '
$ws.Range('D42').Value2 = '0.762'
$ws.Range('E42').Value2 = '  -3.12%  '
$ws.Range('D43').Value2 = '4.31'
$ws.Range('E43').Value2 = '  -2.57%  '
$ws.Range('D44').Value2 = '40.02'
$ws.Range('E44').Value2 = '  -0.30%  '
$ws.Range('D45').Value2 = '0.0656'
$ws.Range('E45').Value2 = '  -2.38%  '
$ws.Range('D46').Value2 = '5.90'
$ws.Range('E46').Value2 = '  -4.23%  '
$ws.Range('D47').Value2 = '2.28'
$ws.Range('E47').Value2 = '  -2.71%  '
$ws.Range('D48').Value2 = '314.53'
$ws.Range('E48').Value2 = '  -1.66%  '
$ws.Range('D49').Value2 = '23.03'
$ws.Range('E49').Value2 = '  -5.12%  '
$ws.Range('E50').Value2 = '  -1.89%  '
$ws.Range('D51').Value2 = '0.104'
$ws.Range('E51').Value2 = '  +3.50%  '

# Drop the temporary Text number format again so untouched formatting
# (style index) matches the rest of the sheet.
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
